# Update "Countries & Provincias Spain" worksheet with refreshed COVID data.
# The underlying data refresh re-sorted provinces 21-31 (descending by total
# cases) after Murcia and Tenerife received updated figures, and a handful
# of Canary Island rows had their "Muertes" (deaths) figure revised from 21
# to 24. The timestamp footer was also bumped forward by 30 minutes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 21-31: new province order + refreshed totals/active/recovered/deaths
$rows = @(
    @{ Row = 21; Name = "Murcia";              B = 596; C = 9;   D = 579; E = 8  },
    @{ Row = 22; Name = "Granada";              B = 579; C = 1;   D = 553; E = 25 },
    @{ Row = 23; Name = "Gipuzkoa/Guipuzcoa";    B = 563; C = 466; D = 543; E = 20 },
    @{ Row = 24; Name = "Tenerife";              B = 539; C = 15;  D = 519; E = 24 },
    @{ Row = 25; Name = "Sevilla";               B = 535; C = 6;   D = 511; E = 18 },
    @{ Row = 26; Name = "Salamanca";             B = 533; C = 57;  D = 430; E = 46 },
    @{ Row = 27; Name = "Cantabria";             B = 510; C = 12;  D = 484; E = 14 },
    @{ Row = 28; Name = "Valladolid";            B = 501; C = 36;  D = 444; E = 21 },
    @{ Row = 29; Name = "Caceres";               B = 485; C = 3;   D = 447; E = 35 },
    @{ Row = 30; Name = "Burgos";                B = 485; C = 55;  D = 402; E = 28 },
    @{ Row = 31; Name = "Leon";                  B = 438; C = 31;  D = 376; E = 31 }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.Name
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
    $ws.Cells.Item($r.Row, 5).Value = $r.E
}

# Canary Island rows whose "Muertes" figure was revised from 21 to 24
$deathUpdateRows = @(45, 55, 57, 58, 62, 63)
foreach ($row in $deathUpdateRows) {
    $ws.Cells.Item($row, 5).Value = 24
}

# Update the "last updated" footer timestamp
$ws.Range("A1").Value = "Datos actualizados a 25 de Marzo de 2020 a las 21:46"
